# 自动更新Excel文件
# For each data row: column D = total days (总天), column E = remaining days (剩余),
# column F = start date (开始时间, yyyyMMdd). The script simulates a daily refresh
# that recomputes "remaining days" against the new reference date. When a cycle
# has fully elapsed (remaining days would hit zero or below), the cycle restarts:
# the start date is reset to the new reference date and the remaining days are
# reset back to the total day count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New reference ("today") date for this automatic update run.
$today = [datetime]::ParseExact("20260215", "yyyyMMdd", $null)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $totalDays = $dCell.Value()
    $startRaw = $fCell.Value()

    if ($null -eq $totalDays -or $null -eq $startRaw) {
        continue
    }

    $startStr = [string]$startRaw

    try {
        $startDate = [datetime]::ParseExact($startStr, "yyyyMMdd", $null)
    } catch {
        # Malformed / unparsable start date -> leave row untouched.
        continue
    }

    $daysElapsed = [int]($today.ToOADate() - $startDate.ToOADate())
    $remaining = [int]$totalDays - $daysElapsed

    if ($remaining -le 0) {
        # Cycle finished - restart it from today.
        $eCell.Value = [int]$totalDays
        $fCell.Value = [int]$today.ToString("yyyyMMdd")
    } else {
        $eCell.Value = $remaining
    }
}
